$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.197.42"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "3.475.31"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").Value = "3.474.41"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("E10").Value = "  -5.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.81%  "
$ws.Range("E12").Value = "  -6.81%  "
$ws.Range("D13").Value = "4.074.76"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000182"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.23%  "
$ws.Range("D16").Value = "3.429.92"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "65.099.86"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.93%  "
$ws.Range("E23").Value = "  -5.58%  "
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.89%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "3.617.79"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.17%  "
$ws.Range("E32").Value = "  -10.12%  "
$ws.Range("D33").Value = "3.492.47"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -6.90%  "
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "171.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0775"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.83%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.67%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.14%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.97%  "
